$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.517.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.84%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'3.486.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'590.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'168.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.73%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'3.483.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.48%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.126"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.16%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.43%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.090.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'28.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.69%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'66.545.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.92%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.483.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "  +3.08%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'13.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.01%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'390.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'7.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'72.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.70%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'0.532"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.20%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "  +5.69%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'10.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.12%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'6.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.08%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "  +4.41%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'23.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'7.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.02%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "  +8.50%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'162.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.30%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "  +3.83%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "  +4.99%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'6.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.52%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.0742"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.48%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'4.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.32%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'26.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'26.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'43.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.764.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "  +2.53%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'344.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.75%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.08"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'33.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +12.47%  "
$ws.Range("E51").Style = "Normal"
